{"js": "// Update the division-problem worksheet table: replace the 25 \"NN\u00f7N=\"\n// problems (in the 5 populated rows of the single table) with their new\n// values, in document/reading order. Row indices below refer to the\n// table's row collection (which includes the blank spacer rows), and\n// column indices are 0-based within each row.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"53\u00f72=\", newText: \"95\u00f78=\" },\n  { row: 0, col: 1, oldText: \"52\u00f75=\", newText: \"42\u00f72=\" },\n  { row: 0, col: 2, oldText: \"16\u00f74=\", newText: \"56\u00f74=\" },\n  { row: 0, col: 3, oldText: \"86\u00f72=\", newText: \"87\u00f77=\" },\n  { row: 0, col: 4, oldText: \"53\u00f73=\", newText: \"25\u00f72=\" },\n\n  { row: 4, col: 0, oldText: \"64\u00f72=\", newText: \"61\u00f78=\" },\n  { row: 4, col: 1, oldText: \"19\u00f73=\", newText: \"38\u00f73=\" },\n  { row: 4, col: 2, oldText: \"10\u00f73=\", newText: \"56\u00f73=\" },\n  { row: 4, col: 3, oldText: \"82\u00f78=\", newText: \"47\u00f76=\" },\n  { row: 4, col: 4, oldText: \"65\u00f77=\", newText: \"48\u00f72=\" },\n\n  { row: 8, col: 0, oldText: \"60\u00f72=\", newText: \"61\u00f72=\" },\n  { row: 8, col: 1, oldText: \"64\u00f77=\", newText: \"98\u00f78=\" },\n  { row: 8, col: 2, oldText: \"26\u00f75=\", newText: \"19\u00f75=\" },\n  { row: 8, col: 3, oldText: \"36\u00f78=\", newText: \"53\u00f78=\" },\n  { row: 8, col: 4, oldText: \"29\u00f74=\", newText: \"43\u00f78=\" },\n\n  { row: 12, col: 0, oldText: \"78\u00f76=\", newText: \"42\u00f75=\" },\n  { row: 12, col: 1, oldText: \"95\u00f78=\", newText: \"19\u00f73=\" },\n  { row: 12, col: 2, oldText: \"80\u00f73=\", newText: \"15\u00f76=\" },\n  { row: 12, col: 3, oldText: \"59\u00f77=\", newText: \"81\u00f75=\" },\n  { row: 12, col: 4, oldText: \"26\u00f74=\", newText: \"57\u00f78=\" },\n\n  { row: 16, col: 0, oldText: \"93\u00f77=\", newText: \"77\u00f76=\" },\n  { row: 16, col: 1, oldText: \"73\u00f73=\", newText: \"12\u00f75=\" },\n  { row: 16, col: 2, oldText: \"58\u00f74=\", newText: \"42\u00f75=\" },\n  { row: 16, col: 3, oldText: \"98\u00f75=\", newText: \"54\u00f73=\" },\n  { row: 16, col: 4, oldText: \"75\u00f75=\", newText: \"35\u00f74=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Grab all the cells we need to touch up front and verify current text\n// before writing, so a structural mismatch fails loudly instead of\n// silently corrupting the wrong cell.\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const expected = replacements[i].oldText;\n  const actual = cells[i].value;\n  if (actual !== expected) {\n    throw new Error(\n      `Cell (${replacements[i].row},${replacements[i].col}) text mismatch: expected \"${expected}\" but found \"${actual}\"`\n    );\n  }\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  cells[i].value = replacements[i].newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: replace the 25 \"NN\u00f7N=\"\n# problems (in the 5 populated rows of the single table) with their new\n# values. Table rows/columns are addressed 1-based, as in the Word\n# object model ($t.Cell(row, col)).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"53\u00f72=\"; New = \"95\u00f78=\" },\n    @{ Row = 1;  Col = 2; Old = \"52\u00f75=\"; New = \"42\u00f72=\" },\n    @{ Row = 1;  Col = 3; Old = \"16\u00f74=\"; New = \"56\u00f74=\" },\n    @{ Row = 1;  Col = 4; Old = \"86\u00f72=\"; New = \"87\u00f77=\" },\n    @{ Row = 1;  Col = 5; Old = \"53\u00f73=\"; New = \"25\u00f72=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"64\u00f72=\"; New = \"61\u00f78=\" },\n    @{ Row = 5;  Col = 2; Old = \"19\u00f73=\"; New = \"38\u00f73=\" },\n    @{ Row = 5;  Col = 3; Old = \"10\u00f73=\"; New = \"56\u00f73=\" },\n    @{ Row = 5;  Col = 4; Old = \"82\u00f78=\"; New = \"47\u00f76=\" },\n    @{ Row = 5;  Col = 5; Old = \"65\u00f77=\"; New = \"48\u00f72=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"60\u00f72=\"; New = \"61\u00f72=\" },\n    @{ Row = 9;  Col = 2; Old = \"64\u00f77=\"; New = \"98\u00f78=\" },\n    @{ Row = 9;  Col = 3; Old = \"26\u00f75=\"; New = \"19\u00f75=\" },\n    @{ Row = 9;  Col = 4; Old = \"36\u00f78=\"; New = \"53\u00f78=\" },\n    @{ Row = 9;  Col = 5; Old = \"29\u00f74=\"; New = \"43\u00f78=\" },\n\n    @{ Row = 13; Col = 1; Old = \"78\u00f76=\"; New = \"42\u00f75=\" },\n    @{ Row = 13; Col = 2; Old = \"95\u00f78=\"; New = \"19\u00f73=\" },\n    @{ Row = 13; Col = 3; Old = \"80\u00f73=\"; New = \"15\u00f76=\" },\n    @{ Row = 13; Col = 4; Old = \"59\u00f77=\"; New = \"81\u00f75=\" },\n    @{ Row = 13; Col = 5; Old = \"26\u00f74=\"; New = \"57\u00f78=\" },\n\n    @{ Row = 17; Col = 1; Old = \"93\u00f77=\"; New = \"77\u00f76=\" },\n    @{ Row = 17; Col = 2; Old = \"73\u00f73=\"; New = \"12\u00f75=\" },\n    @{ Row = 17; Col = 3; Old = \"58\u00f74=\"; New = \"42\u00f75=\" },\n    @{ Row = 17; Col = 4; Old = \"98\u00f75=\"; New = \"54\u00f73=\" },\n    @{ Row = 17; Col = 5; Old = \"75\u00f75=\"; New = \"35\u00f74=\" }\n)\n\nforeach ($rep in $replacements) {\n    $cell = $t.Cell($rep.Row, $rep.Col)\n    $cellRange = $cell.Range\n    # Cell range text includes a trailing end-of-cell mark (and, for the\n    # last cell in a row, an end-of-row mark); strip control characters\n    # before comparing against the expected current value.\n    $current = $cellRange.Text -replace \"[\\a\\r\\n]+$\", \"\"\n    if ($current -ne $rep.Old) {\n        throw \"Cell ($($rep.Row),$($rep.Col)) text mismatch: expected '$($rep.Old)' but found '$current'\"\n    }\n    $cellRange.Text = $rep.New\n}\n"}
